# Apply crypto price/volume updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.239.33'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.589.48'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.61%  '
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  +0.62%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0849'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').Value = '1.812.75'
$ws.Range('E12').Value = '  +0.93%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.04'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.13%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.568.88'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.36'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').Value = '26.235.91'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.42'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '213.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.39%  '
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('E23').Value = '  +1.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.83'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.112'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.19'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  -1.37%  '
$ws.Range('E31').Value = '  +1.18%  '
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('E33').Value = '  -1.02%  '
$ws.Range('D34').Value = '1.336.13'
$ws.Range('E34').Value = '  +4.52%  '
$ws.Range('E35').Value = '  -0.67%  '
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('E37').Value = '  -3.19%  '
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.815'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.74'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.10%  '
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.02'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.72%  '
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.767'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.63%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.724.28'
$ws.Range('E45').Value = '  +0.86%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.80'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.57'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.94%  '
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0973'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.998'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.34%  '
